$wb = $excel.ActiveWorkbook

# New record to append (commit: "[Kadastro App] Yeni kayit eklendi: 2962")
$kayitNo  = "2962"
$tarih    = "2025-09-09"
$birim    = "Erdemli"
$parsel   = "1"
$is       = "3B"
$personel = "SEVİL SARAÇER (Tekniker), EMİNE ALANLI KIRCILI (K.Mühendisi)"

function Add-KayitRow($ws) {
    $lastRow = $ws.UsedRange.Rows.Count
    $newRow = $lastRow + 1

    # Every column in this workbook stores its values as plain text (even
    # record numbers, dates and parcel counts), so force a text number
    # format before writing the values to keep them from being
    # auto-coerced into real numbers/dates.
    $rowRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 6))
    $rowRange.NumberFormat = "@"

    $ws.Cells.Item($newRow, 1).Value = $kayitNo
    $ws.Cells.Item($newRow, 2).Value = $tarih
    $ws.Cells.Item($newRow, 3).Value = $birim
    $ws.Cells.Item($newRow, 4).Value = $parsel
    $ws.Cells.Item($newRow, 5).Value = $is
    $ws.Cells.Item($newRow, 6).Value = $personel
}

# "Kayitlar" is the master records sheet (sheet1.xml / rId1)
$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
Add-KayitRow($wsKayitlar)

# "Erdemli" is the per-birim sheet mirroring records for that birim (sheet8.xml / rId8)
$wsErdemli = $wb.Worksheets.Item("Erdemli")
Add-KayitRow($wsErdemli)
